$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.464.44"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "2.624.21"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.68"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.96"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.546"
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("D9").Value = "2.622.75"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.124"
$ws.Range("E10").Value = "  +6.60%  "
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.22"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.351"
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.00"
$ws.Range("E14").Value = "  +0.00%  "
$ws.Range("D15").Value = "3.102.56"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000183"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").Value = "67.448.13"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").Value = "2.627.68"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.25"
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "364.60"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.59"
$ws.Range("E21").Value = "  -4.39%  "
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.15"
$ws.Range("E23").Value = "  +6.68%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("B25").Value = "Aptos"
$ws.Range("C25").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.00"
$ws.Range("E25").Value = "  -2.95%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "66.13"
$ws.Range("E26").Value = "  -7.70%  "
$ws.Range("D27").Value = "2.762.87"
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "585.22"
$ws.Range("E29").Value = "  -6.00%  "
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.42"
$ws.Range("E31").Value = "  -3.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.92"
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("E33").Value = "  -0.81%  "
$ws.Range("E34").Value = "  -2.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.53"
$ws.Range("E36").Value = "  -3.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.96"
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.18"
$ws.Range("E38").Value = "  +2.09%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.44"
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.371"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.30"
$ws.Range("E41").Value = "  -3.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.83"
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.61"
$ws.Range("E43").Value = "  -1.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.19"
$ws.Range("E44").Value = "  -0.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.36"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "155.83"
$ws.Range("E47").Value = "  -1.18%  "
$ws.Range("D48").Value = "0.0₆0290"
$ws.Range("E48").Value = "  -4.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.73"
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.91"
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.624"
$ws.Range("E51").Value = "  -0.43%  "
